$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix time-of-day component for a handful of early datetime entries
# (09:59:59.999 -> 11:00:00.000, i.e. add 1/24 day)
$ws.Range("A6").Value = 37347.45833333334
$ws.Range("A18").Value = 37712.45833333334
$ws.Range("A30").Value = 38078.45833333334
$ws.Range("A42").Value = 38443.45833333334
$ws.Range("A60").Value = 38991.45833333334

# Update refreshed M2 figures for existing rows 252-254
$ws.Range("C252:F252").Value = 1703600000000
$ws.Range("C253:F253").Value = 1719900000000
$ws.Range("C254:F254").Value = 1749500000000

# Append new row 257 with the latest data point, cloning formatting from row 256
$ws.Range("A256:G256").Copy($ws.Range("A257:G257"))
$ws.Range("A257").Value = 45047.41666666666
$ws.Range("B257").Value = "ECONOMICS:AEM2"
$ws.Range("C257:F257").Value = 1855323000000
$ws.Range("G257").Value = 0
